$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the new "through" date
$ws.Name = "Through 2021-11-10"

# Update the row-13 label text ("November (through 11-09)" -> "... 11-10")
$ws.Range("A13").Value = "November (through 11-10)"

# --- Row 13 (November, partial-month) updated counts/rates ---
$ws.Range("F13").Value = 22
$ws.Range("G13").Value = 0.0435

$ws.Range("I13").Value = 36
$ws.Range("J13").Value = 0.027

$ws.Range("L13").Value = 20
$ws.Range("M13").Value = 0.2

$ws.Range("O13").Value = 15
$ws.Range("P13").Value = 0.1176

$ws.Range("Q13").Value = 2
$ws.Range("R13").Value = 61
$ws.Range("S13").Value = 0.0317

# T13 (2021 arrest_made) and V13 (2021 arrest_rate) are removed entirely
$ws.Range("T13").Clear()
$ws.Range("V13").Clear()

# U13 (2021 no_arrest_made) updated
$ws.Range("U13").Value = 68

# --- Row 14 (Total) updated counts/rates ---
$ws.Range("F14").Value = 456
$ws.Range("G14").Value = 0.1041

$ws.Range("I14").Value = 685
$ws.Range("J14").Value = 0.083

$ws.Range("L14").Value = 569
$ws.Range("M14").Value = 0.1109

$ws.Range("O14").Value = 449
$ws.Range("P14").Value = 0.1002

$ws.Range("Q14").Value = 56
$ws.Range("R14").Value = 1064
$ws.Range("S14").Value = 0.05

$ws.Range("T14").Value = 85
$ws.Range("U14").Value = 1427
$ws.Range("V14").Value = 0.0562
